$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3160.1538
$ws.Range("I138").Value = 3417.0908
$ws.Range("J138").Value = 2971.7334
$ws.Range("K138").Value = 10251.2724
$ws.Range("L138").Value = 8915.200199999999
$ws.Range("M138").Value = -5111.2724
$ws.Range("N138").Value = -19195.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 490120.88
$ws.Range("I32").Value = 651785.5600000001
$ws.Range("J32").Value = 15904.4
$ws.Range("K32").Value = 651785.5600000001
$ws.Range("L32").Value = 15904.4
$ws.Range("M32").Value = -651498.5600000001
$ws.Range("N32").Value = -16478.4
$ws.Range("H41").Value = 1289
$ws.Range("I41").Value = 1289
$ws.Range("K41").Value = 1289
$ws.Range("M41").Value = -875
$ws.Range("H63").Value = 3474
$ws.Range("I63").Value = 3280.3076
$ws.Range("J63").Value = 3622.1177
$ws.Range("K63").Value = 3280.3076
$ws.Range("L63").Value = 3622.1177
$ws.Range("M63").Value = -2594.3076
$ws.Range("N63").Value = -4994.1177
$ws.Range("H66").Value = 3474
$ws.Range("I66").Value = 3280.3076
$ws.Range("J66").Value = 3622.1177
$ws.Range("K66").Value = 16401.538
$ws.Range("L66").Value = 18110.5885
$ws.Range("M66").Value = -12969.538
$ws.Range("N66").Value = -24974.5885
$ws.Range("H132").Value = 7780.5454
$ws.Range("I132").Value = 12453
$ws.Range("J132").Value = 5110.5713
$ws.Range("K132").Value = 37359
$ws.Range("L132").Value = 15331.7139
$ws.Range("M132").Value = -34829
$ws.Range("N132").Value = -20391.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 840.3333
$ws.Range("I94").Value = 822.0769
$ws.Range("J94").Value = 870
$ws.Range("K94").Value = 822.0769
$ws.Range("L94").Value = 870
$ws.Range("M94").Value = -371.0769
$ws.Range("N94").Value = -1772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2641
$ws.Range("I62").Value = 2551.25
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2551.25
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1927.25
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2641
$ws.Range("I65").Value = 2551.25
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 12756.25
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -9636.25
$ws.Range("N65").Value = -21240
$ws.Range("H99").Value = 1164.0834
$ws.Range("I99").Value = 956.9
$ws.Range("K99").Value = 956.9
$ws.Range("M99").Value = 541.1
$ws.Range("H126").Value = 1164.0834
$ws.Range("I126").Value = 956.9
$ws.Range("K126").Value = 2870.7
$ws.Range("M126").Value = -400.6999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1128.6
$ws.Range("I113").Value = 510.375
$ws.Range("J113").Value = 1835.1428
$ws.Range("K113").Value = 1531.125
$ws.Range("L113").Value = 5505.428400000001
$ws.Range("M113").Value = 638.875
$ws.Range("N113").Value = -9845.428400000001
$ws.Range("H114").Value = 1005.4
$ws.Range("I114").Value = 131.75
$ws.Range("J114").Value = 4500
$ws.Range("K114").Value = 395.25
$ws.Range("L114").Value = 13500
$ws.Range("M114").Value = 2858.75
$ws.Range("N114").Value = -20008
$ws.Range("H129").Value = 1035.3636
$ws.Range("I129").Value = 538
$ws.Range("J129").Value = 1449.8334
$ws.Range("K129").Value = 1614
$ws.Range("L129").Value = 4349.5002
$ws.Range("M129").Value = 3386
$ws.Range("N129").Value = -14349.5002
$ws.Range("H131").Value = 1030.5333
$ws.Range("J131").Value = 1117.091
$ws.Range("L131").Value = 3351.273
$ws.Range("N131").Value = -13431.273
$ws.Range("H140").Value = 1526.8718
$ws.Range("I140").Value = 1045.1724
$ws.Range("K140").Value = 3135.5172
$ws.Range("M140").Value = 2044.4828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 41668450
$ws.Range("I122").Value = 62501070
$ws.Range("K122").Value = 187503210
$ws.Range("M122").Value = -187500760
$ws.Range("H132").Value = 2672.5
$ws.Range("I132").Value = 2340.087
$ws.Range("J132").Value = 3522
$ws.Range("K132").Value = 7020.261
$ws.Range("L132").Value = 10566
$ws.Range("M132").Value = -4490.261
$ws.Range("N132").Value = -15626

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1313.375
$ws.Range("I68").Value = 1227.6888
$ws.Range("J68").Value = 1663.909
$ws.Range("K68").Value = 1227.6888
$ws.Range("L68").Value = 1663.909
$ws.Range("M68").Value = -478.6887999999999
$ws.Range("N68").Value = -3161.909
$ws.Range("H71").Value = 1313.375
$ws.Range("I71").Value = 1227.6888
$ws.Range("J71").Value = 1663.909
$ws.Range("K71").Value = 6138.444
$ws.Range("L71").Value = 8319.545
$ws.Range("M71").Value = -2394.444
$ws.Range("N71").Value = -15807.545
$ws.Range("H122").Value = 3233.9333
$ws.Range("I122").Value = 2443.4285
$ws.Range("J122").Value = 3925.625
$ws.Range("K122").Value = 7330.2855
$ws.Range("L122").Value = 11776.875
$ws.Range("M122").Value = -4880.2855
$ws.Range("N122").Value = -16676.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 70019
$ws.Range("J28").Value = 70019
$ws.Range("L28").Value = 70019
$ws.Range("N28").Value = -70715
$ws.Range("H31").Value = 70019
$ws.Range("J31").Value = 70019
$ws.Range("L31").Value = 70019
$ws.Range("N31").Value = -70715
$ws.Range("H62").Value = 28232.5
$ws.Range("I62").Value = 3512.8572
$ws.Range("J62").Value = 62840
$ws.Range("K62").Value = 3512.8572
$ws.Range("L62").Value = 62840
$ws.Range("M62").Value = -2888.8572
$ws.Range("N62").Value = -64088
$ws.Range("H65").Value = 28232.5
$ws.Range("I65").Value = 3512.8572
$ws.Range("J65").Value = 62840
$ws.Range("K65").Value = 17564.286
$ws.Range("L65").Value = 314200
$ws.Range("M65").Value = -14444.286
$ws.Range("N65").Value = -320440
$ws.Range("H81").Value = 4821.5884
$ws.Range("I81").Value = 4212.0713
$ws.Range("J81").Value = 7666
$ws.Range("K81").Value = 8424.142599999999
$ws.Range("L81").Value = 15332
$ws.Range("M81").Value = -7363.142599999999
$ws.Range("N81").Value = -17454
$ws.Range("H84").Value = 4821.5884
$ws.Range("I84").Value = 4212.0713
$ws.Range("J84").Value = 7666
$ws.Range("K84").Value = 42120.713
$ws.Range("L84").Value = 76660
$ws.Range("M84").Value = -36816.713
$ws.Range("N84").Value = -87268
$ws.Range("H126").Value = 1462.7333
$ws.Range("I126").Value = 1252.6
$ws.Range("K126").Value = 3757.8
$ws.Range("M126").Value = -1287.8
